# "scraping better with date" - append two newly scraped arrival rows
# (LO3801 / LO3807, Monday Jan 09, tail numbers SP-LIA / SP-LID) to the
# "Main Data" sheet of the arrivals table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main Data")

# --- Row 20: flight #19 ---
$ws.Range("A20").Value = 19
$ws.Range("B20").Value = "Monday, Jan 09"
$ws.Range("C20").Value = "2:15 PM"
$ws.Range("D20").Value = "LO3801"
$ws.Range("E20").Value = "Warsaw"
$ws.Range("F20").Value = "(WAW)"
$ws.Range("G20").Value = "LOT "
$ws.Range("H20").Value = "E75S"
$ws.Range("I20").Value = "(SP-LIA)"
$ws.Range("J20").Value = "2:14 PM"
$ws.Range("K20").ClearFormats()
$ws.Range("L20").Value = "0 hours, -1 minutes"
$ws.Range("M20").ClearFormats()

# --- Row 21: flight #20 ---
$ws.Range("A21").Value = 20
$ws.Range("B21").Value = "Monday, Jan 09"
$ws.Range("C21").Value = "4:05 PM"
$ws.Range("D21").Value = "LO3807"
$ws.Range("E21").Value = "Warsaw"
$ws.Range("F21").Value = "(WAW)"
$ws.Range("G21").Value = "LOT "
$ws.Range("H21").Value = "E75S"
$ws.Range("I21").Value = "(SP-LID)"
$ws.Range("J21").Value = "4:24 PM"
$ws.Range("K21").ClearFormats()
$ws.Range("L21").Value = "0 hours, 19 minutes"
$ws.Range("M21").ClearFormats()
